# Update '想去人数' (interest count) and '最低票价' (lowest price) figures
# across all four worksheets to match the refreshed data snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 110
$ws.Range("F4").Value = 646
$ws.Range("G7").Value = 30
$ws.Range("F9").Value = 11837
$ws.Range("F10").Value = 200
$ws.Range("F15").Value = 244
$ws.Range("F18").Value = 1210
$ws.Range("G18").Value = 30
$ws.Range("F20").Value = 258
$ws.Range("F21").Value = 745
$ws.Range("F23").Value = 280
$ws.Range("F24").Value = 2909
$ws.Range("F26").Value = 3725
$ws.Range("F27").Value = 3725
$ws.Range("F28").Value = 1082
$ws.Range("F33").Value = 1005
$ws.Range("F34").Value = 43
$ws.Range("F36").Value = 262
$ws.Range("F40").Value = 4053
$ws.Range("F41").Value = 4468
$ws.Range("F42").Value = 5501
$ws.Range("F46").Value = 281
$ws.Range("F47").Value = 71
$ws.Range("F50").Value = 111

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 4164
$ws.Range("F12").Value = 791

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 757
$ws.Range("F3").Value = 426
$ws.Range("F4").Value = 69

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 757
$ws.Range("F3").Value = 426
$ws.Range("F4").Value = 69
$ws.Range("F6").Value = 646
$ws.Range("G9").Value = 30
$ws.Range("F11").Value = 11837
$ws.Range("F16").Value = 244
$ws.Range("F18").Value = 1210
$ws.Range("G18").Value = 30
$ws.Range("F20").Value = 258
$ws.Range("F21").Value = 4164
$ws.Range("F22").Value = 745
$ws.Range("F23").Value = 280
$ws.Range("F25").Value = 3725
$ws.Range("F26").Value = 1082
$ws.Range("F30").Value = 1005
$ws.Range("F31").Value = 43
$ws.Range("F33").Value = 262
$ws.Range("F36").Value = 4468
$ws.Range("F40").Value = 281
$ws.Range("F44").Value = 71
$ws.Range("F50").Value = 111
